# templete_nota_permintaan_barang.docx edit
#
# 1) Paragraphs 1 & 2 ("DINAS SOSIAL PROVINSI JAWA TIMUR" / "SEKRETARIAT"):
#    add single-line spacing (w:spacing after=0 line=240 lineRule=auto) and
#    switch the font to "Sans Serif Collection" (ascii/hAnsi/cs, hint=cs) on
#    both the paragraph mark and the run.
# 2) Paragraph 3 (the blank centered line under "SEKRETARIAT"): drop the
#    center justification so it goes back to default (left) alignment.
# 3) Paragraph 5 ("Nomor: ………./…………./……………"): split the number run into
#    three runs (": ……", "…./", "…………./……………") bracketing the middle one
#    with proofErr gramStart/gramEnd markers, as Word does when you retype
#    part of a run.
#
# Because the Font object's NameBi (=> w:cs) setter — and the engine has no
# settable w:hint property at all — does not propagate to the paragraph
# mark's rPr in this host, the font/spacing changes are applied by replacing
# each paragraph's Range with an equivalent WordprocessingML fragment via
# Range.InsertXML, which lets every attribute (including w:hint="cs") be
# written exactly as Word itself would serialize it.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "DINAS SOSIAL PROVINSI JAWA TIMUR" ---
$p1 = $d.Paragraphs(1)
$p1xml = '<w:p ' + $wNs + ' w14:paraId="7A771492" w14:textId="0258767B" w:rsidR="008D1673" w:rsidRPr="00E85509" w:rsidRDefault="00852FEC" w:rsidP="00852FEC"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Sans Serif Collection" w:hAnsi="Sans Serif Collection" w:cs="Sans Serif Collection" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00E85509"><w:rPr><w:rFonts w:ascii="Sans Serif Collection" w:hAnsi="Sans Serif Collection" w:cs="Sans Serif Collection" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>DINAS SOSIAL PROVINSI JAWA TIMUR</w:t></w:r></w:p>'
$p1.Range.InsertXML($p1xml)

# --- Paragraph 2: "SEKRETARIAT" ---
$p2 = $d.Paragraphs(2)
$p2xml = '<w:p ' + $wNs + ' w14:paraId="59E421AD" w14:textId="014D5107" w:rsidR="00852FEC" w:rsidRDefault="00852FEC" w:rsidP="00E85509"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Sans Serif Collection" w:hAnsi="Sans Serif Collection" w:cs="Sans Serif Collection" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00E85509"><w:rPr><w:rFonts w:ascii="Sans Serif Collection" w:hAnsi="Sans Serif Collection" w:cs="Sans Serif Collection" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>SEKRETARIAT</w:t></w:r></w:p>'
$p2.Range.InsertXML($p2xml)

# --- Paragraph 3: blank line, remove centering ---
$p3 = $d.Paragraphs(3)
$p3xml = '<w:p ' + $wNs + ' w14:paraId="00AC85E9" w14:textId="77777777" w:rsidR="00E85509" w:rsidRPr="00E85509" w:rsidRDefault="00E85509" w:rsidP="00E85509"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$p3.Range.InsertXML($p3xml)

# --- Paragraph 5: "Nomor: ………./…………./……………" -> split number run ---
$p5 = $d.Paragraphs(5)
$p5xml = '<w:p ' + $wNs + ' w14:paraId="2DFBD028" w14:textId="50633288" w:rsidR="00852FEC" w:rsidRPr="00E85509" w:rsidRDefault="00852FEC" w:rsidP="00852FEC"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00E85509"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Nomor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00E85509"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>: ……</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>…./</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>…………./……………</w:t></w:r></w:p>'
$p5.Range.InsertXML($p5xml)
